$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Q17 = 5
$ws.Range("Q17").Value = 5

# O19 = 0
$ws.Range("O19").Value = 0

# O20 = MOD(O19+1,$Q$17)
$ws.Range("O20").Formula = '=MOD(O19+1,$Q$17)'

# O21:O38 = shared formula MOD(prev+1,$Q$17)
$ws.Range("O21:O38").Formula = '=MOD(O20+1,$Q$17)'

$ws.Range("Q16").Select()
